$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A1 text value: "One" -> "OneOneOneOne"
$ws.Range("A1").Value = "OneOneOneOne"

# Apply center horizontal + center vertical alignment to A5
$range5 = $ws.Range("A5")
$range5.HorizontalAlignment = -4108  # xlCenter
$range5.VerticalAlignment = -4108    # xlCenter
